$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data: forma_pagamento / Forma de pagamento
$ws.Range("A15").Value = "forma_pagamento"
$ws.Range("B15").Value = "Forma de pagamento"

# Apply the same formatting as the rest of the data rows (row 14, style index 2)
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to mimic the final state observed in the diff
$ws.Range("B23").Select()
